$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    for ($i = 1; $i -le $sec.Headers.Count; $i++) {
        $h = $sec.Headers.Item($i)
        if ($h.Exists) {
            foreach ($shp in $h.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "BTec_Logo-Orange") {
                    $shp.Name = "image2.jpg"
                }
            }
        }
    }
    for ($i = 1; $i -le $sec.Footers.Count; $i++) {
        $f = $sec.Footers.Item($i)
        if ($f.Exists) {
            foreach ($shp in $f.Range.InlineShapes) {
                if ($shp.AlternativeText -eq "Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png") {
                    $shp.Name = "image1.png"
                }
            }
        }
    }
}

Write-Output "Renamed inline images."
